# Applies the crypto-price/volume refresh described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.775.32'
$ws.Range('E2').Value = '  +0.76%  '

$ws.Range('D3').Value = '1.949.83'
$ws.Range('E3').Value = '  +1.70%  '

$ws.Range('E4').Value = '  -0.15%  '

$ws.Range('D5').Value = '''247.33'
$ws.Range('E5').Value = '  +1.25%  '

$ws.Range('E6').Value = '  -0.10%  '

$ws.Range('D7').Value = '''0.4809'
$ws.Range('E7').Value = '  -1.99%  '

$ws.Range('D8').Value = '''0.2957'
$ws.Range('E8').Value = '  +1.66%  '

$ws.Range('D9').Value = '''0.06824'
$ws.Range('E9').Value = '  +1.71%  '

$ws.Range('D10').Value = '''113.08'
$ws.Range('E10').Value = '  +6.11%  '

$ws.Range('D11').Value = '''19.56'
$ws.Range('E11').Value = '  +4.30%  '

$ws.Range('D12').Value = '1.953.00'
$ws.Range('E12').Value = '  +1.86%  '

$ws.Range('D13').Value = '''5.570'
$ws.Range('E13').Value = '  +5.54%  '

$ws.Range('D14').Value = '''0.07655'
$ws.Range('E14').Value = '  +0.08%  '

$ws.Range('D15').Value = '''0.6923'
$ws.Range('E15').Value = '  +3.98%  '

$ws.Range('D16').Value = '''298.46'
$ws.Range('E16').Value = '  +7.62%  '

$ws.Range('D17').Value = '30.730.27'
$ws.Range('E17').Value = '  +0.66%  '

$ws.Range('D18').Value = '''13.30'
$ws.Range('E18').Value = '  +3.98%  '

$ws.Range('D19').Value = '''5.694'
$ws.Range('E19').Value = '  +3.67%  '

$ws.Range('D20').Value = '''0.000007703'
$ws.Range('E20').Value = '  +2.20%  '

$ws.Range('D21').Value = '2.194.93'
$ws.Range('E21').Value = '  +1.25%  '

$ws.Range('D22').Value = '''0.9999'
$ws.Range('E22').Value = '  -0.10%  '

$ws.Range('E23').Value = '  -0.19%  '

$ws.Range('D24').Value = '''6.586'
$ws.Range('E24').Value = '  +2.75%  '

$ws.Range('D25').Value = '''9.750'
$ws.Range('E25').Value = '  +3.50%  '

$ws.Range('D26').Value = '''168.05'
$ws.Range('E26').Value = '  +2.27%  '

$ws.Range('D27').Value = '''20.48'
$ws.Range('E27').Value = '  +2.73%  '

$ws.Range('D28').Value = '''2.181'
$ws.Range('E28').Value = '  +3.74%  '

$ws.Range('D29').Value = '''0.1089'
$ws.Range('E29').Value = '  +3.67%  '

$ws.Range('D30').Value = '''1.432'
$ws.Range('E30').Value = '  +1.54%  '

$ws.Range('D31').Value = '''4.554'
$ws.Range('E31').Value = '  +12.77%  '

$ws.Range('D32').Value = '''4.436'
$ws.Range('E32').Value = '  +7.04%  '

$ws.Range('D33').Value = '''0.05066'
$ws.Range('E33').Value = '  +1.33%  '

$ws.Range('D34').Value = '''0.7815'
$ws.Range('E34').Value = '  +7.35%  '

$ws.Range('D35').Value = '''1.166'
$ws.Range('E35').Value = '  +2.71%  '

$ws.Range('D36').Value = '''0.02071'
$ws.Range('E36').Value = '  +1.86%  '

$ws.Range('D38').Value = '''2.709'
$ws.Range('E38').Value = '  +1.22%  '

$ws.Range('D39').Value = '''2.043'
$ws.Range('E39').Value = '  +1.69%  '

$ws.Range('D40').Value = '''111.35'
$ws.Range('E40').Value = '  +0.12%  '

$ws.Range('D41').Value = '''0.4480'
$ws.Range('E41').Value = '  +1.09%  '

$ws.Range('D42').Value = '''0.8766'
$ws.Range('E42').Value = '  +0.89%  '

$ws.Range('D43').Value = '''5.980'
$ws.Range('E43').Value = '  +1.69%  '

$ws.Range('D44').Value = '''71.66'
$ws.Range('E44').Value = '  +5.71%  '

$ws.Range('D45').Value = '''1.003'
$ws.Range('E45').Value = '  +0.20%  '

$ws.Range('D46').Value = '''7.424'
$ws.Range('E46').Value = '  +2.46%  '

$ws.Range('D47').Value = '''9.524'
$ws.Range('E47').Value = '  +3.11%  '

$ws.Range('D48').Value = '''49.06'
$ws.Range('E48').Value = '  +1.62%  '

$ws.Range('D49').Value = '''0.1260'
$ws.Range('E49').Value = '  +1.08%  '

# Rows 50 and 51: "WOONetwork" and "Elrond" swap positions (re-ranked)
# and pick up refreshed Price / Volume(1h) figures.
$ws.Range('B50').Value = 'Elrond'
$ws.Range('C50').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('D50').Value = '''35.62'
$ws.Range('E50').Value = '  +2.79%  '

$ws.Range('B51').Value = 'WOONetwork'
$ws.Range('C51').Value = 'https://coinranking.com/coin/k-J3YwacF+woonetwork-woo'
$ws.Range('D51').Value = '''0.2560'
$ws.Range('E51').Value = '  +2.71%  '
